$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts old D->E, old E->F)
$ws.Columns.Item(4).Insert()

# New header for the inserted column
$ws.Range("D1").Value = "Review Count"

# Review Count values for rows 2-22
$reviewCounts = @{
    2 = 2
    3 = 2
    4 = 2
    5 = 2
    6 = 2
    7 = 2
    8 = 2
    9 = 2
    10 = 2
    11 = 4
    12 = 2
    13 = 2
    14 = 2
    15 = 2
    16 = 3
    17 = 2
    18 = 2
    19 = 2
    20 = 2
    21 = 2
    22 = 2
}

foreach ($row in $reviewCounts.Keys) {
    $ws.Cells.Item($row, 4).Value = $reviewCounts[$row]
}
